$wb = $excel.ActiveWorkbook

# --- Add the new "ChangePassword" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ChangePassword"

# --- Header row ---
$ws.Range("A1").Value = "Mật khẩu cũ"
$ws.Range("B1").Value = "Mật khẩu mới"
$ws.Range("C1").Value = "Xác nhận lại mật khẩu"
$ws.Range("D1").Value = "Kết quả mong muốn"

# --- Row 2 ---
$ws.Range("A2").Value = "ngthiquyen"
$ws.Range("B2").Value = "ngthiquyen1"
$ws.Range("C2").Value = "ngthiquyen1"
$ws.Range("D2").Value = "Mật khẩu không đúng"

# --- Row 3 ---
$ws.Range("A3").Value = "ngthiquyen102"
$ws.Range("B3").Value = "ngthiquyen1"
$ws.Range("C3").Value = "ngthiquyen"
$ws.Range("D3").Value = "Xác nhận mật khẩu không khớp"

# --- Row 4 (A4 left blank) ---
$ws.Range("B4").Value = "ngthiquyen2"
$ws.Range("C4").Value = "ngthiquyen2"
$ws.Range("D4").Value = "Vui lòng điền vào trường này."

# --- Row 5 (B5 left blank) ---
$ws.Range("A5").Value = "ngthiquyen102"
$ws.Range("C5").Value = "ngthiquyen3"
$ws.Range("D5").Value = "Vui lòng điền vào trường này."

# --- Row 6 (C6 left blank) ---
$ws.Range("A6").Value = "ngthiquyen102"
$ws.Range("B6").Value = "ngthiquyen"
$ws.Range("D6").Value = "Vui lòng điền vào trường này."

# --- Row 7 / Row 8 were authored out of strict row order: D8's text was
# entered before the rest of row 7, so reproduce that exact sequence to
# keep the shared-string table identical. ---
$ws.Range("A7").Value = "ngthiquyen102"
$ws.Range("D8").Value = "Mật khẩu mới dài từ 6 đến 50 ký tự"

# --- Row 7 (taller row, wrapped text result) ---
$ws.Range("B7").Value = "ad1"
$ws.Range("C7").Value = "ad"
$ws.Range("D7").Value = "Mật khẩu mới dài từ 6 đến 50 ký tự`nXác nhận mật khẩu không khớp"
$ws.Range("D7").Font.Name = "Calibri"
$ws.Range("D7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 29

# --- Row 8 ---
$ws.Range("A8").Value = "ngthiquyen102"
$ws.Range("B8").Value = "as"
$ws.Range("C8").Value = "as"

# --- Row 9 ---
$ws.Range("A9").Value = "ngthiquyen102"
$ws.Range("B9").Value = "ngthiquyen"
$ws.Range("C9").Value = "ngthiquyen"
$ws.Range("D9").Value = "Đổi password thành công"

# --- Column widths (best-fit, matching authored sheet) ---
$ws.Columns.Item(1).ColumnWidth = 13.26953125
$ws.Columns.Item(2).ColumnWidth = 12.453125
$ws.Columns.Item(3).ColumnWidth = 19.1796875
$ws.Columns.Item(4).ColumnWidth = 30.26953125

# --- Select the ChangePassword sheet as active (this is now the last/active tab) ---
$ws.Range("A1").Select()
$ws.Activate()
